$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update example data so it matches the readme
$ws.Range("C5").Value = 3
$ws.Range("C6").Value = 4

# Update the active selection on the sheet
$ws.Range("E19").Select()
